$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency Price (D) and Volume(1h) (E) values.
# Values are prefixed with a leading apostrophe so Excel keeps number-like
# strings (e.g. "242.39") stored as text instead of converting them to a
# floating point number, matching the original inline-string cell content.
# The style is then reset to Normal so no stray quote-prefix style is left
# behind on the cell.

$ws.Range('D2').Value = "'" + '36.481.39'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  -1.41%  '
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').Value = "'" + '2.056.76'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  +0.67%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('E4').Value = "'" + '  -0.20%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').Value = "'" + '242.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  -1.56%  '
$ws.Range('E5').Style = 'Normal'

$ws.Range('E6').Value = "'" + '  +0.57%  '
$ws.Range('E6').Style = 'Normal'

$ws.Range('E7').Value = "'" + '  +0.01%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('D8').Value = "'" + '54.29'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  -5.24%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('D9').Value = "'" + '58.66'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  -2.02%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').Value = "'" + '0.359'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  -5.36%  '
$ws.Range('E10').Style = 'Normal'

$ws.Range('E11').Value = "'" + '  -2.59%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('E12').Value = "'" + '  -3.03%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').Value = "'" + '0.908'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  +2.27%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').Value = "'" + '14.68'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  -5.84%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').Value = "'" + '2.359.25'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  +0.49%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('E16').Value = "'" + '  -5.40%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').Value = "'" + '2.046.06'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  -0.19%  '
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').Value = "'" + '36.407.80'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  -1.50%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').Value = "'" + '16.72'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  -8.49%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').Value = "'" + '71.95'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  -3.04%  '
$ws.Range('E20').Style = 'Normal'

$ws.Range('E21').Value = "'" + '  -4.64%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').Value = "'" + '237.88'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  +0.96%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('E23').Value = "'" + '  -3.50%  '
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').Value = "'" + '0.999'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  -0.10%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('E25').Value = "'" + '  -4.00%  '
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').Value = "'" + '9.30'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  -1.87%  '
$ws.Range('E26').Style = 'Normal'

$ws.Range('E27').Value = "'" + '  -1.26%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('D28').Value = "'" + '164.43'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  -3.28%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('D29').Value = "'" + '20.06'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  +0.68%  '
$ws.Range('E29').Style = 'Normal'

$ws.Range('E30').Value = "'" + '  -1.47%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('D31').Value = "'" + '1.19'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  +5.63%  '
$ws.Range('E31').Style = 'Normal'

$ws.Range('E32').Value = "'" + '  -6.95%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('E33').Value = "'" + '  -5.82%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('E34').Value = "'" + '  -3.50%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('E35').Value = "'" + '  -0.13%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('E36').Value = "'" + '  -0.56%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('E37').Value = "'" + '  -2.86%  '
$ws.Range('E37').Style = 'Normal'

$ws.Range('D38').Value = "'" + '0.0822'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  -5.95%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('E39').Value = "'" + '  -5.81%  '
$ws.Range('E39').Style = 'Normal'

$ws.Range('D40').Value = "'" + '4.86'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  -4.33%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('E41').Value = "'" + '  -3.53%  '
$ws.Range('E41').Style = 'Normal'

$ws.Range('D42').Value = "'" + '2.83'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  -9.10%  '
$ws.Range('E42').Style = 'Normal'

$ws.Range('E43').Value = "'" + '  -3.23%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').Value = "'" + '93.58'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  -3.77%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('D45').Value = "'" + '0.0903'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  -8.43%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('D46').Value = "'" + '1.392.98'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  +7.64%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('D47').Value = "'" + '7.50'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  +10.68%  '
$ws.Range('E47').Style = 'Normal'

$ws.Range('E48').Value = "'" + '  -7.26%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('E49').Value = "'" + '  -0.20%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('E50').Value = "'" + '  -4.26%  '
$ws.Range('E50').Style = 'Normal'

$ws.Range('D51').Value = "'" + '2.247.85'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  +0.78%  '
$ws.Range('E51').Style = 'Normal'
